$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("League Summary")

# Remove the "Team PA", "Team BB", "Team SF" columns (D:F). This shifts the
# former G:J ("Team BA", "Team OBP", "Team SLG", "Team OPS") left into D:G
# and updates the sheet dimension automatically.
$ws.Range("D1:F1").EntireColumn.Delete()

# Resize columns to match the new layout. Excel's ColumnWidth property is
# offset from the stored OOXML column width by 5/6, so subtract that to land
# on the exact target widths.
$offset = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 20 - $offset
$ws.Columns.Item(2).ColumnWidth = 15 - $offset
$ws.Columns.Item(3).ColumnWidth = 16 - $offset
$ws.Columns.Item(4).ColumnWidth = 12 - $offset
$ws.Columns.Item(5).ColumnWidth = 12 - $offset
$ws.Columns.Item(6).ColumnWidth = 12 - $offset
$ws.Columns.Item(7).ColumnWidth = 12 - $offset
